# Fix id error on diagram
# Add an "Out of scope" column to the "Software Component" table (Table1),
# marking MyProcess as "No" (in scope) and every other software component
# as "Yes" (out of scope).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Software Component")
$tbl = $ws.ListObjects.Item("Table1")

# Grow the table by one column (extends ref/autoFilter/tableColumns).
$null = $tbl.ListColumns.Add()

# Header
$ws.Range("D1").Value = "Out of scope"

# Data - row 2 is MyProcess (kept in scope), the rest are out of scope.
$ws.Range("D2").Value = "No"
$ws.Range("D3:D12").Value = "Yes"
